# Updating the RTM according to the new SRS modifications
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------------
# 1. Update CYRS requirement labels (column C) - the old "Req _ DIGELV _CYRS_0x_V1.0"
#    labels are replaced with shorter "DIGELV _CYRS_0x_V1.0" labels, and a brand
#    new row (#9) requirement "DIGELV _CYRS_05_V1.0" is introduced.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = " DIGELV _CYRS_01_V1.0"
$ws.Range("C5").Value = "DIGELV _CYRS_02_V1.0"
$ws.Range("C7").Value = "DIGELV _CYRS_03_V1.0"
$ws.Range("C10").Value = "DIGELV _CYRS_04_V1.0"
$ws.Range("C12").Value = "DIGELV _CYRS_05_V1.0"

# ---------------------------------------------------------------------------
# 2. Update the related SRS traceability lists (column E) to reflect the
#    renumbered / newly added SRS requirements.
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "[DIGELV_SRS_001_V1.0]`n[DIGELV_SRS_002_V1.0]"
$ws.Range("E6").Value = "[DIGELV_SRS_003_V1.0]`n[DIGELV_SRS_004_V1.0]`n[DIGELV_SRS_005_V1.0]`n[DIGELV_SRS_006_V1.0]`n[DIGELV_SRS_007_V1.0]`n[DIGELV_SRS_007_V1.0]`n[DIGELV_SRS_009_V1.0]`n[DIGELV_SRS_010_V1.0]`n[DIGELV_SRS_011_V1.0]"
$ws.Range("E7").Value = "[DIGELV_SRS_012_V1.0]`n[DIGELV_SRS_013_V1.0]"
$ws.Range("E10").Value = "[DIGELV_SRS_014_V1.0]`n[DIGELV_SRS_015_V1.0]"
$ws.Range("E12").Value = "[DIGELV_SRS_016_V1.0]"

# Entering the multi-line text above makes the engine auto-fit rows 4 and 10;
# restore their original (unchanged) height of 15.75 points.
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3. E4 used to be a standalone cell (no horizontal centering); it is now
#    merged together with E5 underneath it. Merge first, then re-apply the
#    per-cell alignment/wrap so E4 keeps the wrapped, centered text while E5
#    remains centered without wrapping (matching the target workbook).
# ---------------------------------------------------------------------------
[void]$ws.Range("E4:E5").Merge()

$e4 = $ws.Range("E4")
$e4.HorizontalAlignment = $xlCenter
$e4.VerticalAlignment = $xlCenter
$e4.WrapText = $true

$e5 = $ws.Range("E5")
$e5.HorizontalAlignment = $xlCenter
$e5.VerticalAlignment = $xlCenter
$e5.WrapText = $false

# ---------------------------------------------------------------------------
# 4. Row 6 needs to grow taller to accommodate the longer SRS list that now
#    lives in E6.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 159

# ---------------------------------------------------------------------------
# 5. Update the active selection saved in the sheet view.
# ---------------------------------------------------------------------------
[void]$ws.Range("B13").Select()
